$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 31282.857
$ws.Range("J87").Value = 31282.857
$ws.Range("L87").Value = 31282.857
$ws.Range("N87").Value = -33778.857
$ws.Range("H90").Value = 31282.857
$ws.Range("J90").Value = 31282.857
$ws.Range("L90").Value = 93848.571
$ws.Range("N90").Value = -106328.571
$ws.Range("H116").Value = 2927.3635
$ws.Range("I116").Value = 2768.8667
$ws.Range("K116").Value = 2768.8667
$ws.Range("M116").Value = 673.1333
$ws.Range("H132").Value = 4253.8887
$ws.Range("I132").Value = 1532.591
$ws.Range("J132").Value = 16227.6
$ws.Range("K132").Value = 4597.772999999999
$ws.Range("L132").Value = 48682.8
$ws.Range("M132").Value = -2067.772999999999
$ws.Range("N132").Value = -53742.8
$ws.Range("H135").Value = 433.65
$ws.Range("I135").Value = 204.05556
$ws.Range("K135").Value = 1836.50004
$ws.Range("M135").Value = 698.4999599999999
$ws.Range("H137").Value = 3000678
$ws.Range("I137").Value = 3517887.8
$ws.Range("J137").Value = 2386491.5
$ws.Range("K137").Value = 10553663.4
$ws.Range("L137").Value = 7159474.5
$ws.Range("M137").Value = -10551113.4
$ws.Range("N137").Value = -7164574.5
$ws.Range("H138").Value = 1916.1072
$ws.Range("J138").Value = 2629.5
$ws.Range("L138").Value = 7888.5
$ws.Range("N138").Value = -18168.5
$ws.Range("H141").Value = 1241
$ws.Range("I141").Value = 431.66666
$ws.Range("J141").Value = 2455
$ws.Range("K141").Value = 1294.99998
$ws.Range("L141").Value = 7365
$ws.Range("M141").Value = 3885.00002
$ws.Range("N141").Value = -17725

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 353684.03
$ws.Range("I61").Value = 251273.34
$ws.Range("J61").Value = 594650.3
$ws.Range("K61").Value = 251273.34
$ws.Range("L61").Value = 594650.3
$ws.Range("M61").Value = -251061.34
$ws.Range("N61").Value = -595074.3
$ws.Range("H63").Value = 2633.3333
$ws.Range("I63").Value = 2450
$ws.Range("K63").Value = 2450
$ws.Range("M63").Value = -1764
$ws.Range("H66").Value = 2633.3333
$ws.Range("I66").Value = 2450
$ws.Range("K66").Value = 12250
$ws.Range("M66").Value = -8818
$ws.Range("H102").Value = 12600
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622
$ws.Range("H122").Value = 3135.1333
$ws.Range("I122").Value = 2097.6365
$ws.Range("J122").Value = 5988.25
$ws.Range("K122").Value = 6292.9095
$ws.Range("L122").Value = 17964.75
$ws.Range("M122").Value = -3842.9095
$ws.Range("N122").Value = -22864.75
$ws.Range("H136").Value = 353684.03
$ws.Range("I136").Value = 251273.34
$ws.Range("J136").Value = 594650.3
$ws.Range("K136").Value = 753820.02
$ws.Range("L136").Value = 1783950.9
$ws.Range("M136").Value = -751270.02
$ws.Range("N136").Value = -1789050.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1130.93
$ws.Range("I31").Value = 630.86957
$ws.Range("J31").Value = 2243.9678
$ws.Range("K31").Value = 630.86957
$ws.Range("L31").Value = 2243.9678
$ws.Range("M31").Value = -335.86957
$ws.Range("N31").Value = -2833.9678
$ws.Range("H34").Value = 1130.93
$ws.Range("I34").Value = 630.86957
$ws.Range("J34").Value = 2243.9678
$ws.Range("K34").Value = 630.86957
$ws.Range("L34").Value = 2243.9678
$ws.Range("M34").Value = -428.86957
$ws.Range("N34").Value = -2647.9678
$ws.Range("H58").Value = 3478.182
$ws.Range("I58").Value = 4067.3667
$ws.Range("J58").Value = 2215.6428
$ws.Range("K58").Value = 4067.3667
$ws.Range("L58").Value = 2215.6428
$ws.Range("M58").Value = -3864.3667
$ws.Range("N58").Value = -2621.6428
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H132").Value = 13516282
$ws.Range("I132").Value = 20835434
$ws.Range("J132").Value = 4000.923
$ws.Range("K132").Value = 62506302
$ws.Range("L132").Value = 12002.769
$ws.Range("M132").Value = -62503772
$ws.Range("N132").Value = -17062.769
$ws.Range("H134").Value = 11365465
$ws.Range("I134").Value = 14707298
$ws.Range("J134").Value = 3230.2
$ws.Range("K134").Value = 44121894
$ws.Range("L134").Value = 9690.599999999999
$ws.Range("M134").Value = -44119359
$ws.Range("N134").Value = -14760.6
$ws.Range("H136").Value = 3478.182
$ws.Range("I136").Value = 4067.3667
$ws.Range("J136").Value = 2215.6428
$ws.Range("K136").Value = 12202.1001
$ws.Range("L136").Value = 6646.928400000001
$ws.Range("M136").Value = -9652.1001
$ws.Range("N136").Value = -11746.9284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 317.85715
$ws.Range("I4").Value = 153.83333
$ws.Range("J4").Value = 1302
$ws.Range("K4").Value = 461.49999
$ws.Range("L4").Value = 3906
$ws.Range("M4").Value = -349.49999
$ws.Range("N4").Value = -4130
$ws.Range("H6").Value = 58823990
$ws.Range("I6").Value = 166666930
$ws.Range("J6").Value = 574.1818
$ws.Range("K6").Value = 500000790
$ws.Range("L6").Value = 1722.5454
$ws.Range("M6").Value = -500000677
$ws.Range("N6").Value = -1948.5454
$ws.Range("H9").Value = 105002340
$ws.Range("J9").Value = 157502500
$ws.Range("L9").Value = 472507500
$ws.Range("N9").Value = -472507948
$ws.Range("H10").Value = 206.16667
$ws.Range("I10").Value = 52.57143
$ws.Range("J10").Value = 421.2
$ws.Range("K10").Value = 157.71429
$ws.Range("L10").Value = 1263.6
$ws.Range("M10").Value = -18.71429000000001
$ws.Range("N10").Value = -1541.6
$ws.Range("H11").Value = 200.33333
$ws.Range("I11").Value = 150.5
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 451.5
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = -311.5
$ws.Range("N11").Value = -1180
$ws.Range("H13").Value = 650.125
$ws.Range("J13").Value = 728.7143
$ws.Range("L13").Value = 2186.1429
$ws.Range("N13").Value = -2522.1429
$ws.Range("H15").Value = 414.2143
$ws.Range("I15").Value = 147.75
$ws.Range("J15").Value = 520.8
$ws.Range("K15").Value = 443.25
$ws.Range("L15").Value = 1562.4
$ws.Range("M15").Value = -303.25
$ws.Range("N15").Value = -1842.4
$ws.Range("H16").Value = 750
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3346
$ws.Range("H19").Value = 500
$ws.Range("J19").Value = 500
$ws.Range("L19").Value = 1500
$ws.Range("N19").Value = -1848
$ws.Range("H68").Value = 1100.5862
$ws.Range("J68").Value = 1266.775
$ws.Range("L68").Value = 3800.325
$ws.Range("N68").Value = -5422.325000000001
$ws.Range("H71").Value = 1100.5862
$ws.Range("J71").Value = 1266.775
$ws.Range("L71").Value = 11400.975
$ws.Range("N71").Value = -19512.975
$ws.Range("H107").Value = 1229.5
$ws.Range("I107").Value = 350
$ws.Range("J107").Value = 2109
$ws.Range("K107").Value = 1050
$ws.Range("L107").Value = 6327
$ws.Range("M107").Value = 870
$ws.Range("N107").Value = -10167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 6550
$ws.Range("I21").Value = 6200
$ws.Range("J21").Value = 6666.6665
$ws.Range("K21").Value = 6200
$ws.Range("L21").Value = 6666.6665
$ws.Range("M21").Value = -6026
$ws.Range("N21").Value = -7014.6665
$ws.Range("H24").Value = 8141.2
$ws.Range("I24").Value = 9206
$ws.Range("J24").Value = 7875
$ws.Range("K24").Value = 9206
$ws.Range("L24").Value = 7875
$ws.Range("M24").Value = -8863
$ws.Range("N24").Value = -8561
$ws.Range("H132").Value = 13344561
$ws.Range("I132").Value = 4399.3335
$ws.Range("J132").Value = 25658556
$ws.Range("K132").Value = 13198.0005
$ws.Range("L132").Value = 76975668
$ws.Range("M132").Value = -10668.0005
$ws.Range("N132").Value = -76980728
$ws.Range("H136").Value = 3508.2334
$ws.Range("I136").Value = 1961.8788
$ws.Range("J136").Value = 5398.222
$ws.Range("K136").Value = 5885.636399999999
$ws.Range("L136").Value = 16194.666
$ws.Range("M136").Value = -3335.636399999999
$ws.Range("N136").Value = -21294.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2495.4546
$ws.Range("I81").Value = 850
$ws.Range("J81").Value = 5375
$ws.Range("K81").Value = 1700
$ws.Range("L81").Value = 10750
$ws.Range("M81").Value = -639
$ws.Range("N81").Value = -12872
$ws.Range("H84").Value = 2495.4546
$ws.Range("I84").Value = 850
$ws.Range("J84").Value = 5375
$ws.Range("K84").Value = 8500
$ws.Range("L84").Value = 53750
$ws.Range("M84").Value = -3196
$ws.Range("N84").Value = -64358
